$d = $word.ActiveDocument

# --- 1) Title paragraph: "Test3" -> "Service Writer Check in Procedure", centered, Arial Rounded MT Bold 22pt
$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.Find.Execute("Test3", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Service Writer Check in Procedure", 2)

$titlePara.Alignment = 1
$titleRange2 = $titlePara.Range
$titleRange2.Font.Name = "Arial Rounded MT Bold"
$titleRange2.Font.Size = 22
$titleRange2.Font.SizeBi = 22

# remove the stray _GoBack bookmark that used to trail the title text
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- 2) "Add mileage in..." bullet: change trailing phrase
$d.Content.Find.Execute( `
    "Add mileage in and add inspection for CBS service check requiring technicians to check all fluids and service needs", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Add mileage in and add inspection for Complimentary Safety Check", 2)

# --- 3) "Have customer sign tablet..." bullet: add "approve via text or"
$d.Content.Find.Execute( `
    "Have customer sign tablet for authorization of work to be performed", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Have customer approve via text or sign tablet for authorization of work to be performed", 2)

# --- 4) Remove the ""Clean up" revisions..." bullet entirely
$cleanupPara = $d.Paragraphs.Item(10)
$cleanupPara.Range.Delete()
